# Added filtering options for the Component Analysis
# Clears out specific cells in rows 2,3,5,6,7 that should no longer
# contain values (mirrors removing columns from the component table).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2:K2").ClearContents()
$ws.Range("I3:K3").ClearContents()
$ws.Range("K5").ClearContents()
$ws.Range("J6:K6").ClearContents()
$ws.Range("I7:K7").ClearContents()
